$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.058134008642892
$ws.Range("D2").Value = 1.06245570516235
$ws.Range("E2").Value = 1.053963137347903
$ws.Range("F2").Value = 1.07120195294803
$ws.Range("I2").Value = 1.047987611823276
$ws.Range("J2").Value = 1.063126810246374
$ws.Range("K2").Value = 1.065177004793279
$ws.Range("L2").Value = 1.056707657155978
$ws.Range("M2").Value = 1.073899755022631

# Row 3
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.059862937996661
$ws.Range("D3").Value = 1.063832423060787
$ws.Range("E3").Value = 1.055481907089025
$ws.Range("F3").Value = 1.072715604502015
$ws.Range("I3").Value = 1.048502448399858
$ws.Range("J3").Value = 1.064504742034891
$ws.Range("K3").Value = 1.066367129738078
$ws.Range("L3").Value = 1.058037771461454
$ws.Range("M3").Value = 1.075228193111517

# Row 4
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.060978667832776
$ws.Range("D4").Value = 1.064720581980623
$ws.Range("E4").Value = 1.056461518941101
$ws.Range("F4").Value = 1.073692580818814
$ws.Range("I4").Value = 1.048832764465003
$ws.Range("J4").Value = 1.065393040585174
$ws.Range("K4").Value = 1.067134021105519
$ws.Range("L4").Value = 1.058894816529603
$ws.Range("M4").Value = 1.076084818791368

# Row 5
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.061447016811283
$ws.Range("D5").Value = 1.065093335774026
$ws.Range("E5").Value = 1.056872610763556
$ws.Range("F5").Value = 1.074102725050797
$ws.Range("I5").Value = 1.048970959877549
$ws.Range("J5").Value = 1.065765699724993
$ws.Range("K5").Value = 1.067455666623571
$ws.Range("L5").Value = 1.05925426208524
$ws.Range("M5").Value = 1.076444245917837

# Row 6
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.061525613827976
$ws.Range("D6").Value = 1.065155886251876
$ws.Range("E6").Value = 1.056941592053395
$ws.Range("F6").Value = 1.074171556597445
$ws.Range("I6").Value = 1.048994124384609
$ws.Range("J6").Value = 1.065828225353813
$ws.Range("K6").Value = 1.067509628307764
$ws.Range("L6").Value = 1.059314564724135
$ws.Range("M6").Value = 1.07650455468072

# Row 7
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.060984928679716
$ws.Range("D7").Value = 1.064725565188528
$ws.Range("E7").Value = 1.056467014847938
$ws.Range("F7").Value = 1.073698063435706
$ws.Range("I7").Value = 1.048834613662381
$ws.Range("J7").Value = 1.065398023129205
$ws.Range("K7").Value = 1.067138321903632
$ws.Range("L7").Value = 1.058899622806071
$ws.Range("M7").Value = 1.076089624203061

# Row 8
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.058718939341728
$ws.Range("D8").Value = 1.062921532772319
$ws.Range("E8").Value = 1.054477069542585
$ws.Range("F8").Value = 1.071714014488179
$ws.Range("I8").Value = 1.048162189738352
$ws.Range("J8").Value = 1.06359318258187
$ws.Range("K8").Value = 1.065579881671759
$ws.Range("L8").Value = 1.05715793412514
$ws.Range("M8").Value = 1.074349326934912

# Row 9
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.054702296365187
$ws.Range("D9").Value = 1.059721659471584
$ws.Range("E9").Value = 1.050945935208678
$ws.Range("F9").Value = 1.068198517476199
$ws.Range("I9").Value = 1.046955483554464
$ws.Range("J9").Value = 1.060386882351599
$ws.Range("K9").Value = 1.062808747583994
$ws.Range("L9").Value = 1.054060516049459
$ws.Range("M9").Value = 1.071259517233206

# Row 10
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.05200765754475
$ws.Range("D10").Value = 1.057573641600843
$ws.Range("E10").Value = 1.048574471348003
$ws.Range("F10").Value = 1.065841100597635
$ws.Range("I10").Value = 1.04613602571547
$ws.Range("J10").Value = 1.058231137369114
$ws.Range("K10").Value = 1.060943890426831
$ws.Range("L10").Value = 1.051975736321836
$ws.Range("M10").Value = 1.069183362327299

# Row 11
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.050836632389136
$ws.Range("D11").Value = 1.056639867987463
$ws.Range("E11").Value = 1.047543293001744
$ws.Range("F11").Value = 1.064816884663703
$ws.Range("I11").Value = 1.045777563910667
$ws.Range("J11").Value = 1.057293183857205
$ws.Range("K11").Value = 1.060132105555583
$ws.Range("L11").Value = 1.051068124609694
$ws.Range("M11").Value = 1.068280346769588

# Row 12
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.050401007657313
$ws.Range("D12").Value = 1.056292458381118
$ws.Range("E12").Value = 1.047159602349565
$ws.Range("F12").Value = 1.064435914150283
$ws.Range("I12").Value = 1.045643863340229
$ws.Range("J12").Value = 1.056944094849259
$ws.Range("K12").Value = 1.059829915389218
$ws.Range("L12").Value = 1.050730247642171
$ws.Range("M12").Value = 1.067944308074644

# Row 13
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.050494480528784
$ws.Range("D13").Value = 1.056367004684392
$ws.Range("E13").Value = 1.047241935656909
$ws.Range("F13").Value = 1.064517657886601
$ws.Range("I13").Value = 1.045672567643766
$ws.Range("J13").Value = 1.057019007152127
$ws.Range("K13").Value = 1.059894766147056
$ws.Range("L13").Value = 1.050802757563159
$ws.Range("M13").Value = 1.068016417723265

# Row 14
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.050800636959163
$ws.Range("D14").Value = 1.056611162595421
$ws.Range("E14").Value = 1.047511590675233
$ws.Range("F14").Value = 1.064785404441023
$ws.Range("I14").Value = 1.045766523475363
$ws.Range("J14").Value = 1.057264342240736
$ws.Range("K14").Value = 1.060107139901145
$ws.Range("L14").Value = 1.051040210973336
$ws.Range("M14").Value = 1.068252582402512

# Row 15
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.05098918284115
$ws.Range("D15").Value = 1.056761521133536
$ws.Range("E15").Value = 1.047677645410293
$ws.Range("F15").Value = 1.06495030111684
$ws.Range("I15").Value = 1.045824339441111
$ws.Range("J15").Value = 1.057415409229303
$ws.Range("K15").Value = 1.06023790292334
$ws.Range("L15").Value = 1.051186414075332
$ws.Range("M15").Value = 1.068398008889645

# Row 16
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.052085283956068
$ws.Range("D16").Value = 1.057635534659153
$ws.Range("E16").Value = 1.048642814883365
$ws.Range("F16").Value = 1.065909000752556
$ws.Range("I16").Value = 1.046159738599958
$ws.Range("J16").Value = 1.058293290134111
$ws.Range("K16").Value = 1.060997674462921
$ws.Range("L16").Value = 1.052035867241295
$ws.Range("M16").Value = 1.069243206504865

# Row 17
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.052771694054881
$ws.Range("D17").Value = 1.058182788839464
$ws.Range("E17").Value = 1.049247072507124
$ws.Range("F17").Value = 1.066509436867206
$ws.Range("I17").Value = 1.046369149006053
$ws.Range("J17").Value = 1.058842745875995
$ws.Range("K17").Value = 1.06147310177123
$ws.Range("L17").Value = 1.052567387383255
$ws.Range("M17").Value = 1.069772289052647

# Row 18
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.053171658699875
$ws.Range("D18").Value = 1.058501639983027
$ws.Range("E18").Value = 1.049599110156376
$ws.Range("F18").Value = 1.066859330279041
$ws.Range("I18").Value = 1.046490944632909
$ws.Range("J18").Value = 1.059162800785321
$ws.Range("K18").Value = 1.061749997474784
$ws.Range("L18").Value = 1.052876943468189
$ws.Range("M18").Value = 1.070080506347617

# Row 19
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.053307967872827
$ws.Range("D19").Value = 1.058610300520121
$ws.Range("E19").Value = 1.049719075893094
$ws.Range("F19").Value = 1.066978579213264
$ws.Range("I19").Value = 1.046532414632706
$ws.Range("J19").Value = 1.059271858188926
$ws.Range("K19").Value = 1.061844342083202
$ws.Range("L19").Value = 1.052982414723737
$ws.Range("M19").Value = 1.070185535142713

# Row 20
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.052698090953403
$ws.Range("D20").Value = 1.058124110250784
$ws.Range("E20").Value = 1.049182284487371
$ws.Range("F20").Value = 1.066445050076733
$ws.Range("I20").Value = 1.046346717491067
$ws.Range("J20").Value = 1.058783839404237
$ws.Range("K20").Value = 1.061422135730704
$ws.Range("L20").Value = 1.052510409090182
$ws.Range("M20").Value = 1.069715563691799

# Row 21
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.050710499720431
$ws.Range("D21").Value = 1.056539279868396
$ws.Range("E21").Value = 1.047432202484726
$ws.Range("F21").Value = 1.064706574532444
$ws.Range("I21").Value = 1.045738871111203
$ws.Range("J21").Value = 1.057192116376704
$ws.Range("K21").Value = 1.060044619366204
$ws.Range("L21").Value = 1.050970307676625
$ws.Range("M21").Value = 1.0681830549741

# Row 22
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.049457029924785
$ws.Range("D22").Value = 1.055539561229509
$ws.Range("E22").Value = 1.046328000218803
$ws.Range("F22").Value = 1.063610445817474
$ws.Range("I22").Value = 1.045353497663065
$ws.Range("J22").Value = 1.056187329673079
$ws.Range("K22").Value = 1.05917471186443
$ws.Range("L22").Value = 1.049997640403432
$ws.Range("M22").Value = 1.06721592058968

# Row 23
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.050121883685965
$ws.Range("D23").Value = 1.056069845768154
$ws.Range("E23").Value = 1.046913729730893
$ws.Range("F23").Value = 1.064191821222833
$ws.Range("I23").Value = 1.045558096534513
$ws.Range("J23").Value = 1.056720371004387
$ws.Range("K23").Value = 1.05963623154626
$ws.Range("L23").Value = 1.050513686628481
$ws.Range("M23").Value = 1.067728961232715

# Row 24
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.052731350283308
$ws.Range("D24").Value = 1.058150625669728
$ws.Range("E24").Value = 1.04921156068868
$ws.Range("F24").Value = 1.066474144717927
$ws.Range("I24").Value = 1.046356854408303
$ws.Range("J24").Value = 1.058810458041982
$ws.Range("K24").Value = 1.061445166366619
$ws.Range("L24").Value = 1.052536156582181
$ws.Range("M24").Value = 1.069741196638187

# Row 25
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.055743598831268
$ws.Range("D25").Value = 1.060551455367222
$ws.Range("E25").Value = 1.051861818507638
$ws.Range("F25").Value = 1.069109726669906
$ws.Range("I25").Value = 1.047270063626385
$ws.Range("J25").Value = 1.061218940711613
$ws.Range("K25").Value = 1.063528176960803
$ws.Range("L25").Value = 1.054864711110115
$ws.Range("M25").Value = 1.072061125483577
